$d = $word.ActiveDocument

# Locate the very end of the document body content (just before the
# final paragraph mark / sectPr) and insert the new OOXML paragraphs
# there, preserving exact formatting (including the several
# intentionally-empty <w:p/> separator paragraphs) as captured in the
# source revision.
$endPos = $d.Content.End - 1
$rng = $d.Range($endPos, $endPos)

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>What is the benefit of TypeScript?</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>TypeScript adds types to JavaScript to help spot potential bugs in, clarify the structure of, and help refactor our code</w:t></w:r></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>Given the following code, what will happen after running the TS compiler?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>Compiles without any issues.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>TypeScript’s </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>tsc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t> command will let you know if your code tried to access members of variables known to not exist. TypeScript helps us quickly locate bugs in our code.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:t>How many errors will TypeScript show based on the code provided?</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="AppleSystemUIFont" w:hAnsi="AppleSystemUIFont" w:cs="AppleSystemUIFont"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:t>2</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="5"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:br/></w:r><w:r><w:br/></w:r></w:p>
'@

$rng.InsertXML($xml)
Write-Output "Inserted quiz questions 2-4 at end of document"
